$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text update: "Time" -> "Time (hrs)"
$ws.Range("B1").Value = "Time (hrs)"

# Row 8: move the time value out of C8, put a plain hours total in B8 instead,
# leaving C8 present but empty (still time-formatted).
$ws.Range("C8").ClearContents()
$ws.Range("B8").Value = 2

# New row 9: another date + hours entry, with C9/D9 left as empty time-formatted cells.
$ws.Range("A9").NumberFormat = "MM/DD/YY"
$ws.Range("A9").Value = 42937
$ws.Range("B9").Value = 0.66
$ws.Range("C9").NumberFormat = 'HH:MM:SS\ AM/PM'
$ws.Range("D9").NumberFormat = 'HH:MM:SS\ AM/PM'

# New row 10: a leftover time-only entry in C10.
$ws.Range("C10").NumberFormat = 'HH:MM:SS\ AM/PM'
$ws.Range("C10").Value = 0.0763888888888889

# Columns C and D were resized (e.g. to better fit the new values).
$ws.Columns.Item(3).ColumnWidth = 10.75
$ws.Columns.Item(4).ColumnWidth = 9.7

# Selection moved to C10.
$ws.Range("C10").Select()
